$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3884013333333334
$ws.Range("H2").Value = 1.165204
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.227883333333333
$ws.Range("N2").Value = 3.68365
$ws.Range("O2").Value = 0.05271132222573729
$ws.Range("P2").Value = 0.05917739429803119
$ws.Range("Q2").Value = 0.4769115238444445
$ws.Range("R2").Value = 4.2922037146
$ws.Range("S2").Value = 0.05271132222573729
$ws.Range("T2").Value = 0.05917739429803119

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3884013333333334
$ws.Range("H3").Value = 1.165204
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.377739666666666
$ws.Range("N3").Value = 19.133219
$ws.Range("O3").Value = 0.2737874857612962
$ws.Range("P3").Value = 0.3073728625014814
$ws.Range("Q3").Value = 2.477122590186222
$ws.Range("R3").Value = 22.294103311676
$ws.Range("S3").Value = 0.2737874857612962
$ws.Range("T3").Value = 0.3073728625014814

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3884013333333334
$ws.Range("H4").Value = 1.165204
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.839059333333334
$ws.Range("N4").Value = 14.517178
$ws.Range("O4").Value = 0.2077340809703377
$ws.Range("P4").Value = 0.2332167189067104
$ws.Range("Q4").Value = 1.879497097145778
$ws.Range("R4").Value = 16.915473874312
$ws.Range("S4").Value = 0.2077340809703377
$ws.Range("T4").Value = 0.2332167189067104

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3884013333333334
$ws.Range("H5").Value = 1.165204
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.213925999999999
$ws.Range("N5").Value = 9.641777999999999
$ws.Range("O5").Value = 0.1379693692362262
$ws.Range("P5").Value = 0.1548940041643702
$ws.Range("Q5").Value = 1.248293143634666
$ws.Range("R5").Value = 11.234638292712
$ws.Range("S5").Value = 0.1379693692362262
$ws.Range("T5").Value = 0.1548940041643702

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.3884013333333334
$ws.Range("H6").Value = 1.165204
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 7.635880999999999
$ws.Range("N6").Value = 15.271762
$ws.Range("O6").Value = 0.3277977418064026
$ws.Range("P6").Value = 0.2453390201294068
$ws.Range("Q6").Value = 2.965786361574667
$ws.Range("R6").Value = 17.794718169448
$ws.Range("S6").Value = 0.3277977418064026
$ws.Range("T6").Value = 0.2453390201294068
